# Applies the "changed list of Use cases" edit to the Use case list workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text changes -------------------------------------------------------

# WIND.UC.004 "Modify Customer Accounts" -> "Changing Customer Password"
$ws.Range("B5").Value = "Changing Customer Password"

# "What is done?" column (D): fill in the missing "+" marks
$ws.Range("D4").Value = "+"
$ws.Range("D5").Value = "+"

# WIND.UC.009 "What is done?" note: " + |-" -> " + "
$ws.Range("D10").Value = " + "

# Service-instance related use cases renamed (SI -> Service Instance)
$ws.Range("B15").Value = "Creating Service Instance"
$ws.Range("D15").Value = "+"

$ws.Range("B16").Value = "Modifying Parameters for Service Instance"
$ws.Range("D16").Value = "+"

$ws.Range("B17").Value = "Disconnect for Existing Service Instance"
$ws.Range("D17").Value = "+"

# Remaining rows in the "What is done?" column that were blank before
$ws.Range("D18").Value = "+"
$ws.Range("D19").Value = "+"
$ws.Range("D20").Value = "+"
$ws.Range("D21").Value = "+"
$ws.Range("D22").Value = "+"

# --- Highlight WIND.UC.007 / WIND.UC.008 rows in red --------------------

$ws.Range("A8:D9").Font.Color = 255

# --- Column width / selection cosmetics ----------------------------------

$ws.Columns.Item(2).ColumnWidth = 40.25

$ws.Activate()
$ws.Range("G8").Select()
